$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.564.81'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.623.73'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.60'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.22'
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = '1.853.25'
$ws.Range("E12").Value = '  -1.40%  '
$ws.Range("D13").Value = '1.628.32'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.24'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '27.529.33'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.54'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.53'
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.43'
$ws.Range("E22").Value = '  +2.78%  '
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("E24").Value = '  +5.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.67'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").Value = '1.465.12'
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("E34").Value = '  -2.63%  '
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.35'
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  +6.60%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.874'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0167'
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  -3.01%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.52'
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.21'
$ws.Range("E44").Value = '  -2.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.28'
$ws.Range("E45").Value = '  -6.05%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.763.49'
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.75'
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.43'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.71'
$ws.Range("E51").Value = '  -1.35%  '
